# Add a new "readme_file" row to the "study" sheet of the schema template,
# right before the existing "summary_statistics_assembly" row (new row 18),
# pushing the rows below it down by one. The new field is optional
# everywhere (not mandatory for metadata, not mandatory for summary stats).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("study")

# Insert a fresh row at position 18 (shifts old rows 18-20 down to 19-21),
# inheriting formatting from the surrounding rows.
$ws.Rows.Item(18).Insert()

# Fill in the new row. Columns (row 1 headers): A=NAME, B=DESCRIPTION,
# C=DEFAULT, D=MANDATORY, E=MANDATORY-METADATA, F=MANDATORY-SUMMARY_STATS,
# G=TYPE, L=HEADER.
# NB: assign the HEADER (L) value before the longer DESCRIPTION (B) value so
# the newly created shared-string table entries land in the same order as
# the reference document (readme_file, Readme file, <description>).
$ws.Cells.Item(18, 1).Value = "readme_file"
$ws.Cells.Item(18, 12).Value = "Readme file"
$ws.Cells.Item(18, 2).Value = "Path or URL pointing to the file containing additional information for the provided summary statistics file"
$ws.Cells.Item(18, 3).Value = $true
$ws.Cells.Item(18, 4).Value = $false
$ws.Cells.Item(18, 5).Value = $false
$ws.Cells.Item(18, 6).Value = $false
$ws.Cells.Item(18, 7).Value = "string"

# The hidden _FilterDatabase name for "study" needs to grow by one row to
# keep covering the autofilter header column through the new last row.
$fd = $wb.Names.Item("study!_FilterDatabase")
$fd.RefersTo = "=study!`$C`$1:`$C`$19"

# Make "study" the active sheet/tab and park the selection on the newly
# added row's second column, matching where the edit was made.
$ws.Activate()
$ws.Range("B19").Select()
